$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = -7.877599999999997
$ws.Range("C7").Value = -13.1014
$ws.Range("A9").Value = -21.86310000000001
$ws.Range("C12").Value = -10.6962
$ws.Range("A13").Value = -22.3047
$ws.Range("C14").Value = -13.90749999999999
$ws.Range("D15").Value = -9.071599999999989
$ws.Range("A16").Value = -21.56089999999999
$ws.Range("A18").Value = -22.23040000000001
$ws.Range("C19").Value = -11.9271
$ws.Range("A20").Value = -21.73099999999998
$ws.Range("A26").Value = -21.71929999999999
$ws.Range("C26").Value = -13.0626
$ws.Range("A27").Value = -21.80109999999998
$ws.Range("C27").Value = -12.70539999999999
$ws.Range("D28").Value = -8.411199999999997
$ws.Range("A29").Value = -20.59139999999998
$ws.Range("C29").Value = -11.1576
$ws.Range("D33").Value = -7.844299999999996
$ws.Range("A35").Value = -20.1772
$ws.Range("D35").Value = -8.048899999999998
$ws.Range("A36").Value = -20.8884
$ws.Range("C37").Value = -13.3722
$ws.Range("C38").Value = -13.09119999999999
$ws.Range("D38").Value = -8.649599999999994
$ws.Range("D43").Value = -8.415199999999995
$ws.Range("D44").Value = -7.940600000000001
$ws.Range("A45").Value = -21.56979999999998
$ws.Range("D45").Value = -7.755999999999998
$ws.Range("C47").Value = -11.6998
$ws.Range("D47").Value = -7.5455
$ws.Range("C51").Value = -13.1581
$ws.Range("D51").Value = -7.999899999999994
$ws.Range("C52").Value = -11.2436
$ws.Range("D54").Value = -8.0793
$ws.Range("A55").Value = -22.16900000000001
$ws.Range("C55").Value = -13.02789999999999
$ws.Range("A57").Value = -21.9668
$ws.Range("D57").Value = -8.245799999999999
$ws.Range("D62").Value = -8.352299999999996
$ws.Range("D63").Value = -8.0101
$ws.Range("D67").Value = -6.093
$ws.Range("A69").Value = -21.7024
$ws.Range("C69").Value = -10.9744
$ws.Range("C70").Value = -12.75600000000001
$ws.Range("D70").Value = -8.0786
$ws.Range("A76").Value = -19.59079999999999
$ws.Range("C76").Value = -12.82820000000001
$ws.Range("A78").Value = -19.86329999999998
$ws.Range("C81").Value = -13.5016
$ws.Range("D81").Value = -8.154800000000002
$ws.Range("A82").Value = -22.1607
$ws.Range("A83").Value = -21.94729999999999
$ws.Range("C83").Value = -13.91809999999999
$ws.Range("D88").Value = -7.991199999999996
$ws.Range("A93").Value = -20.53179999999999
$ws.Range("C94").Value = -10.3558
$ws.Range("D96").Value = -7.785000000000004
$ws.Range("A97").Value = -21.8526
$ws.Range("D99").Value = -7.902299999999997
$ws.Range("C100").Value = -12.85989999999999
$ws.Range("C102").Value = -13.0053